# Update the two date cells on the active sheet and move the current
# selection to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# O2: 44803 (2022-08-30) -> 44774 (2022-08-01)
$ws.Range("O2").Value = 44774

# P2: 44803 (2022-08-30) -> 44798 (2022-08-25)
$ws.Range("P2").Value = 44798

# Move the active selection from T2 to P2
$ws.Range("P2").Select() | Out-Null
